$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 9, shifting the existing rows
# 9..62 down to 10..63 (carrying their formatting along, matching the
# "insert row" semantics Excel uses when a new weekly data point is
# recorded at the top of the series).
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new weekly record.
$ws.Range("A9").Value = 11
$ws.Range("B9").Value = "Vega Monumental Concepción"
$ws.Range("C9").Value = "Bíobío"
$ws.Range("D9").Value = 44881
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 100112031
$ws.Range("G9").Value = "Poroto verde"
$ws.Range("H9").Value = "Magnum"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 45000
$ws.Range("L9").Value = 46000
$ws.Range("M9").Value = 45500
$ws.Range("N9").Value = "$/saco 25 kilos"
$ws.Range("O9").Value = "Región de O'Higgins"
$ws.Range("P9").Value = 1820
$ws.Range("Q9").Value = 25
$ws.Range("R9").Value = "Hortaliza"
